$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("A1").Value = 'Vaisselle '
$ws.Range("B1").Value = '''Vaisselle du quotidien : garde le nombre nécessaire pour tenir entre deux vaisselles<br>Vaisselle pour les occasions : en as-tu réellement besoin ? Un service que tu utilises une fois par an ne sert à rien.<br>Cassé/abîmé : à jeter ou recycler.'
$ws.Range("B1").WrapText = $true
$ws.Range("C1").Value = 'et verres, couverts, tasses, bols…'
$ws.Rows.Item(1).RowHeight = 28.8

$ws.Range("A2").Value = 'Casseroles'
$ws.Range("B2").Value = '''Garde le matériel de cuisson qui est en bon état : les revêtements abîmés dégagent des particules potentiellement toxiques dans la nourriture !'
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = 'poêles, chaudrons, marmittes, faitouts…'
$ws.Rows.Item(2).RowHeight = 28.8

$ws.Range("A3").Value = 'Ustensiles'
$ws.Range("B3").Value = '''Élimine les doublons et va au plus simple en gardant des ustensiles multi-usages :<br>- une fourchette comme presse-purée<br>- des bols comme saucière ou pour les apéro<br>- tes mains pour presser un citron...<br> Tout ce que tu n''utilises pas tous les mois au minimum est probablement superflu.'
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = 'Objectif : plus rien sur le plan de travail, tout dans les tiroirs !'
$ws.Rows.Item(3).RowHeight = 43.2

$ws.Range("A4").Value = 'Linge de table'
$ws.Range("B4").Value = 'Débarrasse-toi du linge qui n''est pas pratique (mauvaise taille, manque d''absorption…) ou abîmé. Tu peux égalemen troquer tes produits jetables (essuie-tout, éponge…) contre des produits lavables, à choisir dans des couleurs faciles d''entretien.'
$ws.Range("B4").WrapText = $true
$ws.Range("C4").Value = 'Assez de nappes pour recevoir à Versailles'
$ws.Rows.Item(4).RowHeight = 28.8

$ws.Range("A5").Value = 'Boîtes de conservation'
$ws.Range("B5").Value = 'Garde les boîtes qui ont la bonne taille, qui ne sont pas déformées ou usées et dont tu possèdes le couvercle. Range-les avec le couvercle pour éviter de devoir chercher à chaque utilisation. Conserve le nombre de boîtes nécessaire selon ton utilisation.<br>Astuce : un bocal de récupération peut servir à transporter ou congeler des aliments.'
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = 'Il est où ce couvercle ???'
$ws.Rows.Item(5).RowHeight = 43.2

$ws.Range("A6").Value = 'Gadgets'
$ws.Range("B6").Value = 'Débarrasse-toi de tous les gadgets achetés sur un coup de tête et dont tu ne te sers jamais.'
$ws.Range("B6").WrapText = $true
$ws.Range("C6").Value = 'Ils étaient pourtant censés te simplifier la vie…'

$ws.Range("A7").Value = 'Électroménager'
$ws.Range("B7").Value = 'Tout ce qui ne sert qu''une fois ou deux par an peut être emprunté ! Privilégie également les appareils multi-usages : un grill qui permet de faire des croque-monsieur, des paninis, de la viande, des plancha…, un appareil de cuisson avec lequel tu peux réaliser plusieurs préparations différentes...'
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value = 'Le comptoir de ma cuisine Cap Canaveral'
$ws.Rows.Item(7).RowHeight = 43.2

$ws.Range("A8").Value = 'Livres de cuisine et recettes'
$ws.Range("B8").Value = 'Note, scanne ou photographie les recettes dont tu as vraiment besoin, ne garde pas de livres de recettes complets si tu ne les utilise (quasiment) jamais.<br><a href="https://www.theflonicles.be/2018/09/jai-simplifie-mes-repas-bien-manger.html">Simplifie tes repas</a> pour moins te prendre la tête.'
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = 'Pas le meilleur bouquin que j''ai lu…'
$ws.Rows.Item(8).RowHeight = 43.2

$ws.Range("A9").Value = 'Garde-manger et congélateur'
$ws.Range("B9").Value = 'Fais régulièrement le tri dans tes stocks de nourriture : les aliments périment et s''altèrent avec le temps. Évite d''acheter de grandes quantités d''avance et fais l''inventaire régulièrement, ça ne sert à rien de remplir son congélateur et de garnir toutes les étagères du garde-manger.'
$ws.Range("B9").WrapText = $true
$ws.Range("C9").Value = 'Prêt.e à affronter la prochaine catastrophe mondiale !'
$ws.Rows.Item(9).RowHeight = 43.2

$ws.Range("A10").Value = 'Produits ménagers'
$ws.Range("B10").Value = 'Tu peux nettoyer toute la maison avec très peu de produits ménagers (un savon de Marseille toutes surfaces par exemple). Si tu veux te débarrasser de produits ménagers superflus ou qui ne te conviennent plus, de nombreuses associations (refuges pour animaux par exemple) les récupèreront avec plaisir !<br>Tu peux acheter <a href="https://www.theflonicles.be/2018/12/comme-avant-cosmetiques-naturels-zerodechet.html>un pain de savon ou des paillettes à dissoudre dans l''eau</a> pour doser à ta guise, c''est une solution écologique et économique !'
$ws.Range("B10").WrapText = $true
$ws.Range("C10").Value = 'Nettoyer, balayer, astiquer… (c''est cadeau)'
$ws.Rows.Item(10).RowHeight = 72

$ws.Columns.Item(2).ColumnWidth = 105.21875
$ws.Columns.Item(3).ColumnWidth = 50.6640625

$ws.Range("B9").Select() | Out-Null

Write-Output "done"